# Update scripts with new TPM values.
# Original row 2 (FAPs -> Adm2 -> Calcr -> MuSCs) is recomputed with new TPM
# data and split into two rows: one for target cluster "ECs" and one for the
# (recomputed) "MuSCs" target cluster.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: FAPs -> Adm2 -> Calcr -> ECs (new target cluster, replacing the old
# MuSCs row's numbers with recomputed TPM-based values)
$ws.Range("D2").Value = "ECs"
$ws.Range("H2").Value = 0.5797639999999999
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.09584066666666667
$ws.Range("N2").Value = 0.287522
$ws.Range("O2").Value = 0.3176649051884244
$ws.Range("P2").Value = 0.3176649051884244
$ws.Range("Q2").Value = 0.01852165608977778
$ws.Range("R2").Value = 0.166694904808
$ws.Range("S2").Value = 0.3176649051884244
$ws.Range("T2").Value = 0.3176649051884244

# Row 3 (new): FAPs -> Adm2 -> Calcr -> MuSCs, with recomputed TPM values
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Adm2"
$ws.Range("C3").Value = "Calcr"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.1932546666666667
$ws.Range("H3").Value = 0.5797639999999999
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.205863
$ws.Range("N3").Value = 0.6175889999999999
$ws.Range("O3").Value = 0.6823350948115756
$ws.Range("P3").Value = 0.6823350948115755
$ws.Range("Q3").Value = 0.039783985444
$ws.Range("R3").Value = 0.3580558689959999
$ws.Range("S3").Value = 0.6823350948115756
$ws.Range("T3").Value = 0.6823350948115755
